$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$names = @(
    'Construction Materials(8)',
    'Multiline Retail(7)',
    'Multi-Utilities(18)',
    'Containers & Packaging(12)',
    'Machinery(86)',
    'Road & Rail(22)',
    'Auto Components(21)',
    'Air Freight & Logistics(11)',
    'Electric Utilities(28)',
    'Airlines(14)',
    'Trading Companies & Distributors(25)',
    'Gas Utilities(12)',
    'Equity Real Estate Investment Trusts ...(98)',
    'Hotels, Restaurants & Leisure(51)',
    'Specialty Retail(59)',
    'Wireless Telecommunication Services(14)',
    'Media(42)',
    'Chemicals(52)',
    'Construction & Engineering(21)',
    'Leisure Products(11)',
    'Semiconductors & Semiconductor Equipment(70)',
    'Household Durables(39)',
    'Household Products(10)',
    'Marine(15)',
    'Aerospace & Defense(37)',
    'Building Products(24)',
    'Capital Markets(76)',
    'Communications Equipment(45)',
    'Oil, Gas & Consumable Fuels(125)',
    'Energy Equipment & Services(38)',
    'Insurance(75)',
    'Technology Hardware, Storage & Periph...(19)',
    'Commercial Services & Supplies(52)',
    'Consumer Finance(15)',
    'Textiles, Apparel & Luxury Goods(29)',
    'Diversified Consumer Services(17)',
    'Water Utilities(13)',
    'Diversified Telecommunication Services(20)',
    'ETF(303)',
    'Professional Services(35)',
    'Electrical Equipment(28)',
    'Life Sciences Tools & Services(19)',
    'Electronic Equipment, Instruments & C...(78)',
    'Internet & Direct Marketing Retail(15)',
    'Software(70)',
    'IT Services(52)',
    'Health Care Providers & Services(47)',
    'Banks(251)',
    'Food & Staples Retailing(16)',
    'Beverages(21)',
    'Personal Products(19)',
    'Health Care Equipment & Supplies(86)',
    'Metals & Mining(106)',
    'Thrifts & Mortgage Finance(47)',
    'Entertainment(22)',
    'Food Products(46)',
    'Biotechnology(128)',
    'Pharmaceuticals(53)'
)

$vals = @(
    0.6668117385209243,
    0.637158856603528,
    0.5909373411224138,
    0.5820552138888447,
    0.574297532791705,
    0.5710044321418385,
    0.5537157241117897,
    0.5118791143994464,
    0.5001381885300904,
    0.4917499218936111,
    0.4854130637614995,
    0.4786730779818642,
    0.4718664674272164,
    0.47141191301275,
    0.4711957565821278,
    0.4702411239861925,
    0.4632595642572921,
    0.4613830175459043,
    0.45996107255264,
    0.4273113347199632,
    0.418425023374555,
    0.4173637377507899,
    0.4064015974987705,
    0.4033073324332358,
    0.3975633353773633,
    0.3962955128265593,
    0.3908099843454869,
    0.3803307251461573,
    0.3799339577785524,
    0.3737127086805898,
    0.3678833024531595,
    0.3647382387362086,
    0.3607899127967236,
    0.3521445803575914,
    0.3515114243143465,
    0.3499089347939527,
    0.336959255027089,
    0.3291281176173813,
    0.3252695302660542,
    0.3241847688660761,
    0.3177397643706062,
    0.3115688160483179,
    0.3103076408935375,
    0.2952069581098644,
    0.2918107675410762,
    0.2895560816946511,
    0.2880824409626052,
    0.2879755813017612,
    0.2841416381203662,
    0.270311358287111,
    0.2646998600836321,
    0.2526533277781496,
    0.2491614927465716,
    0.2454875116583609,
    0.2426693261706819,
    0.1892125073366443,
    0.1561710254254096,
    0.154122602581187
)

for ($i = 0; $i -lt $names.Count; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 1).Value = $names[$i]
    $ws.Cells.Item($row, 2).Value = $vals[$i]
}
